$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7169949999999999
$ws.Range("H2").Value = 2.150985
$ws.Range("I2").Value = 0.003651663653539308
$ws.Range("J2").Value = 0.003651663653539308
$ws.Range("O2").Value = 0.8215990550008899
$ws.Range("P2").Value = 0.82159905500089
$ws.Range("Q2").Value = 0.06067450688333333
$ws.Range("R2").Value = 0.5460705619499999
$ws.Range("S2").Value = 0.003000203406928993
$ws.Range("T2").Value = 0.003000203406928993
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7169949999999999
$ws.Range("H3").Value = 2.150985
$ws.Range("I3").Value = 0.003651663653539308
$ws.Range("J3").Value = 0.003651663653539308
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.018375
$ws.Range("N3").Value = 0.05512499999999999
$ws.Range("O3").Value = 0.17840094499911
$ws.Range("P3").Value = 0.17840094499911
$ws.Range("Q3").Value = 0.013174783125
$ws.Range("R3").Value = 0.118573048125
$ws.Range("S3").Value = 0.0006514602466103153
$ws.Range("T3").Value = 0.0006514602466103153
$ws.Range("I4").Value = 0.9751961860217362
$ws.Range("J4").Value = 0.9751961860217361
$ws.Range("O4").Value = 0.8215990550008899
$ws.Range("P4").Value = 0.82159905500089
$ws.Range("S4").Value = 0.8012202648759306
$ws.Range("T4").Value = 0.8012202648759306
$ws.Range("I5").Value = 0.9751961860217362
$ws.Range("J5").Value = 0.9751961860217361
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.018375
$ws.Range("N5").Value = 0.05512499999999999
$ws.Range("O5").Value = 0.17840094499911
$ws.Range("P5").Value = 0.17840094499911
$ws.Range("Q5").Value = 3.518395852999999
$ws.Range("R5").Value = 31.66556267699999
$ws.Range("S5").Value = 0.1739759211458056
$ws.Range("T5").Value = 0.1739759211458056
$ws.Range("G6").Value = 4.138615666666666
$ws.Range("H6").Value = 12.415847
$ws.Range("I6").Value = 0.02107801645190694
$ws.Range("J6").Value = 0.02107801645190694
$ws.Range("O6").Value = 0.8215990550008899
$ws.Range("P6").Value = 0.82159905500089
$ws.Range("Q6").Value = 0.3502234530988889
$ws.Range("R6").Value = 3.15201107789
$ws.Range("S6").Value = 0.01731767839817996
$ws.Range("T6").Value = 0.01731767839817995
$ws.Range("G7").Value = 4.138615666666666
$ws.Range("H7").Value = 12.415847
$ws.Range("I7").Value = 0.02107801645190694
$ws.Range("J7").Value = 0.02107801645190694
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.018375
$ws.Range("N7").Value = 0.05512499999999999
$ws.Range("O7").Value = 0.17840094499911
$ws.Range("P7").Value = 0.17840094499911
$ws.Range("Q7").Value = 0.07604706287499999
$ws.Range("R7").Value = 0.6844235658749999
$ws.Range("S7").Value = 0.003760338053726987
$ws.Range("T7").Value = 0.003760338053726986
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.014556
$ws.Range("H8").Value = 0.043668
$ws.Range("I8").Value = 0.00007413387281768795
$ws.Range("J8").Value = 0.00007413387281768795
$ws.Range("O8").Value = 0.8215990550008899
$ws.Range("P8").Value = 0.82159905500089
$ws.Range("Q8").Value = 0.00123177724
$ws.Range("R8").Value = 0.01108599516
$ws.Range("S8").Value = 0.00006090831985056858
$ws.Range("T8").Value = 0.00006090831985056858
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.014556
$ws.Range("H9").Value = 0.043668
$ws.Range("I9").Value = 0.00007413387281768795
$ws.Range("J9").Value = 0.00007413387281768795
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.018375
$ws.Range("N9").Value = 0.05512499999999999
$ws.Range("O9").Value = 0.17840094499911
$ws.Range("P9").Value = 0.17840094499911
$ws.Range("Q9").Value = 0.0002674665
$ws.Range("R9").Value = 0.002407198499999999
$ws.Range("S9").Value = 0.00001322555296711936
$ws.Range("T9").Value = 0.00001322555296711936

Write-Output "applied 98 cell changes"